$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Daniel 1:1'
$ws.Cells.Item(3, 1).Value = 'Daniel 1:2'
$ws.Cells.Item(4, 1).Value = 'Daniel 1:3'
$ws.Cells.Item(5, 1).Value = 'Daniel 1:4'
$ws.Cells.Item(6, 1).Value = 'Daniel 1:5'
$ws.Cells.Item(7, 1).Value = 'Daniel 1:6'
$ws.Cells.Item(8, 1).Value = 'Daniel 1:7'
$ws.Cells.Item(9, 1).Value = 'Daniel 1:8'
$ws.Cells.Item(10, 1).Value = 'Daniel 1:10'
$ws.Cells.Item(11, 1).Value = 'Daniel 1:11'
$ws.Cells.Item(12, 1).Value = 'Daniel 1:12'
$ws.Cells.Item(13, 1).Value = 'Daniel 1:13'
$ws.Cells.Item(14, 1).Value = 'Daniel 1:14'
$ws.Cells.Item(15, 1).Value = 'Daniel 1:15'
$ws.Cells.Item(16, 1).Value = 'Daniel 1:17'
$ws.Cells.Item(17, 1).Value = 'Daniel 1:18'
$ws.Cells.Item(18, 1).Value = 'Daniel 1:19'
$ws.Cells.Item(19, 1).Value = 'Daniel 1:20'
$ws.Cells.Item(20, 1).Value = 'Daniel 1:21'
$ws.Cells.Item(21, 1).Value = 'Daniel 2:1'
$ws.Cells.Item(22, 1).Value = 'Daniel 2:2'
$ws.Cells.Item(23, 1).Value = 'Daniel 2:3'
$ws.Cells.Item(24, 1).Value = 'Daniel 2:4'
$ws.Cells.Item(25, 1).Value = 'Daniel 2:5'
$ws.Cells.Item(26, 1).Value = 'Daniel 2:6'
$ws.Cells.Item(27, 1).Value = 'Daniel 2:7'
$ws.Cells.Item(28, 1).Value = 'Daniel 2:8'
$ws.Cells.Item(30, 1).Value = 'Daniel 2:11'
$ws.Cells.Item(31, 1).Value = 'Daniel 2:12'
$ws.Cells.Item(32, 1).Value = 'Daniel 2:13'
$ws.Cells.Item(33, 1).Value = 'Daniel 2:14'
$ws.Cells.Item(34, 1).Value = 'Daniel 2:15'
$ws.Cells.Item(35, 1).Value = 'Daniel 2:16'
$ws.Cells.Item(36, 1).Value = 'Daniel 2:17'
$ws.Cells.Item(37, 1).Value = 'Daniel 2:19'
$ws.Cells.Item(38, 1).Value = 'Daniel 2:20'
$ws.Cells.Item(39, 1).Value = 'Daniel 2:21'
$ws.Cells.Item(40, 1).Value = 'Daniel 2:22'
$ws.Cells.Item(41, 1).Value = 'Daniel 2:23'
$ws.Cells.Item(42, 1).Value = 'Daniel 2:24'
$ws.Cells.Item(43, 1).Value = 'Daniel 2:25'
$ws.Cells.Item(44, 1).Value = 'Daniel 2:26'
$ws.Cells.Item(45, 1).Value = 'Daniel 2:27'
$ws.Cells.Item(46, 1).Value = 'Daniel 2:28'
$ws.Cells.Item(47, 1).Value = 'Daniel 2:29'
$ws.Cells.Item(48, 1).Value = 'Daniel 2:30'
$ws.Cells.Item(49, 1).Value = 'Daniel 2:31'
$ws.Cells.Item(50, 1).Value = 'Daniel 2:34'
$ws.Cells.Item(51, 1).Value = 'Daniel 2:35'
$ws.Cells.Item(52, 1).Value = 'Daniel 2:36'
$ws.Cells.Item(53, 1).Value = 'Daniel 2:37'
$ws.Cells.Item(54, 1).Value = 'Daniel 2:38'
$ws.Cells.Item(55, 1).Value = 'Daniel 2:40'
$ws.Cells.Item(56, 1).Value = 'Daniel 2:41'
$ws.Cells.Item(57, 1).Value = 'Daniel 2:42'
$ws.Cells.Item(58, 1).Value = 'Daniel 2:43'
$ws.Cells.Item(59, 1).Value = 'Daniel 2:44'
$ws.Cells.Item(60, 1).Value = 'Daniel 2:45'
$ws.Cells.Item(61, 1).Value = 'Daniel 2:46'
$ws.Cells.Item(62, 1).Value = 'Daniel 2:47'
$ws.Cells.Item(63, 1).Value = 'Daniel 2:48'
$ws.Cells.Item(64, 1).Value = 'Daniel 2:49'
$ws.Cells.Item(65, 1).Value = 'Daniel 3:1'
$ws.Cells.Item(66, 1).Value = 'Daniel 3:2'
$ws.Cells.Item(67, 1).Value = 'Daniel 3:3'
$ws.Cells.Item(68, 1).Value = 'Daniel 3:4'
$ws.Cells.Item(69, 1).Value = 'Daniel 3:5'
$ws.Cells.Item(70, 1).Value = 'Daniel 3:6'
$ws.Cells.Item(71, 1).Value = 'Daniel 3:7'
$ws.Cells.Item(72, 1).Value = 'Daniel 3:8'
$ws.Cells.Item(73, 1).Value = 'Daniel 3:9'
$ws.Cells.Item(74, 1).Value = 'Daniel 3:10'
$ws.Cells.Item(75, 1).Value = 'Daniel 3:11'
$ws.Cells.Item(76, 1).Value = 'Daniel 3:12'
$ws.Cells.Item(77, 1).Value = 'Daniel 3:13'
$ws.Cells.Item(78, 1).Value = 'Daniel 3:14'
$ws.Cells.Item(79, 1).Value = 'Daniel 3:15'
$ws.Cells.Item(80, 1).Value = 'Daniel 3:17'
$ws.Cells.Item(81, 1).Value = 'Daniel 3:18'
$ws.Cells.Item(82, 1).Value = 'Daniel 3:19'
$ws.Cells.Item(83, 1).Value = 'Daniel 3:20'
$ws.Cells.Item(84, 1).Value = 'Daniel 3:21'
$ws.Cells.Item(85, 1).Value = 'Daniel 3:22'
$ws.Cells.Item(86, 1).Value = 'Daniel 3:23'
$ws.Cells.Item(87, 1).Value = 'Daniel 3:24'
$ws.Cells.Item(88, 1).Value = 'Daniel 3:25'
$ws.Cells.Item(89, 1).Value = 'Daniel 3:26'
$ws.Cells.Item(90, 1).Value = 'Daniel 3:27'
$ws.Cells.Item(91, 1).Value = 'Daniel 3:28'
$ws.Cells.Item(92, 1).Value = 'Daniel 3:29'
$ws.Cells.Item(93, 1).Value = 'Daniel 3:30'
$ws.Cells.Item(94, 1).Value = 'Daniel 4:1'
$ws.Cells.Item(95, 1).Value = 'Daniel 4:2'
$ws.Cells.Item(96, 1).Value = 'Daniel 4:3'
$ws.Cells.Item(97, 1).Value = 'Daniel 4:5'
$ws.Cells.Item(98, 1).Value = 'Daniel 4:6'
$ws.Cells.Item(99, 1).Value = 'Daniel 4:7'
$ws.Cells.Item(100, 1).Value = 'Daniel 4:8'
$ws.Cells.Item(101, 1).Value = 'Daniel 4:9'
$ws.Cells.Item(102, 1).Value = 'Daniel 4:13'
$ws.Cells.Item(103, 1).Value = 'Daniel 4:14'
$ws.Cells.Item(104, 1).Value = 'Daniel 4:15'
$ws.Cells.Item(105, 1).Value = 'Daniel 4:16'
$ws.Cells.Item(106, 1).Value = 'Daniel 4:17'
$ws.Cells.Item(107, 1).Value = 'Daniel 4:18'
$ws.Cells.Item(108, 1).Value = 'Daniel 4:19'
$ws.Cells.Item(110, 1).Value = 'Daniel 4:23'
$ws.Cells.Item(111, 1).Value = 'Daniel 4:24'
$ws.Cells.Item(112, 1).Value = 'Daniel 4:26'
$ws.Cells.Item(113, 1).Value = 'Daniel 4:27'
$ws.Cells.Item(114, 1).Value = 'Daniel 4:28'
$ws.Cells.Item(115, 1).Value = 'Daniel 4:29'
$ws.Cells.Item(116, 1).Value = 'Daniel 4:30'
$ws.Cells.Item(117, 1).Value = 'Daniel 4:33'
$ws.Cells.Item(118, 1).Value = 'Daniel 4:35'
$ws.Cells.Item(119, 1).Value = 'Daniel 4:36'
$ws.Cells.Item(120, 1).Value = 'Daniel 4:37'
$ws.Cells.Item(121, 1).Value = 'Daniel 5:3'
$ws.Cells.Item(122, 1).Value = 'Daniel 5:4'
$ws.Cells.Item(123, 1).Value = 'Daniel 5:8'
$ws.Cells.Item(124, 1).Value = 'Daniel 5:9'
$ws.Cells.Item(127, 1).Value = 'Daniel 5:13'
$ws.Cells.Item(128, 1).Value = 'Daniel 5:14'
$ws.Cells.Item(129, 1).Value = 'Daniel 5:15'
$ws.Cells.Item(130, 1).Value = 'Daniel 5:17'
$ws.Cells.Item(131, 1).Value = 'Daniel 5:18'
$ws.Cells.Item(132, 1).Value = 'Daniel 5:19'
$ws.Cells.Item(133, 1).Value = 'Daniel 5:20'
$ws.Cells.Item(134, 1).Value = 'Daniel 5:22'
$ws.Cells.Item(135, 1).Value = 'Daniel 5:23'
$ws.Cells.Item(136, 1).Value = 'Daniel 5:24'
$ws.Cells.Item(137, 1).Value = 'Daniel 5:25'
$ws.Cells.Item(138, 1).Value = 'Daniel 5:26'
$ws.Cells.Item(139, 1).Value = 'Daniel 5:29'
$ws.Cells.Item(140, 1).Value = 'Daniel 6:1'
$ws.Cells.Item(141, 1).Value = 'Daniel 6:2'
$ws.Cells.Item(142, 1).Value = 'Daniel 6:3'
$ws.Cells.Item(143, 1).Value = 'Daniel 6:4'
$ws.Cells.Item(144, 1).Value = 'Daniel 6:5'
$ws.Cells.Item(145, 1).Value = 'Daniel 6:6'
$ws.Cells.Item(146, 1).Value = 'Daniel 6:7'
$ws.Cells.Item(147, 1).Value = 'Daniel 6:8'
$ws.Cells.Item(148, 1).Value = 'Daniel 6:9'
$ws.Cells.Item(149, 1).Value = 'Daniel 6:10'
$ws.Cells.Item(150, 1).Value = 'Daniel 6:11'
$ws.Cells.Item(151, 1).Value = 'Daniel 6:12'
$ws.Cells.Item(152, 1).Value = 'Daniel 6:13'
$ws.Cells.Item(153, 1).Value = 'Daniel 6:14'
$ws.Cells.Item(154, 1).Value = 'Daniel 6:15'
$ws.Cells.Item(155, 1).Value = 'Daniel 6:16'
$ws.Cells.Item(156, 1).Value = 'Daniel 6:18'
$ws.Cells.Item(157, 1).Value = 'Daniel 6:19'
$ws.Cells.Item(158, 1).Value = 'Daniel 6:20'
$ws.Cells.Item(159, 1).Value = 'Daniel 6:21'
$ws.Cells.Item(160, 1).Value = 'Daniel 6:22'
$ws.Cells.Item(161, 1).Value = 'Daniel 6:23'
$ws.Cells.Item(162, 1).Value = 'Daniel 6:24'
$ws.Cells.Item(164, 1).Value = 'Daniel 6:27'
$ws.Cells.Item(165, 1).Value = 'Daniel 6:28'
$ws.Cells.Item(166, 1).Value = 'Daniel 7:1'
$ws.Cells.Item(167, 1).Value = 'Daniel 7:2'
$ws.Cells.Item(168, 1).Value = 'Daniel 7:3'
$ws.Cells.Item(169, 1).Value = 'Daniel 7:4'
$ws.Cells.Item(170, 1).Value = 'Daniel 7:5'
$ws.Cells.Item(171, 1).Value = 'Daniel 7:6'
$ws.Cells.Item(172, 1).Value = 'Daniel 7:7'
$ws.Cells.Item(173, 1).Value = 'Daniel 7:8'
$ws.Cells.Item(174, 1).Value = 'Daniel 7:9'
$ws.Cells.Item(175, 1).Value = 'Daniel 7:10'
$ws.Cells.Item(176, 1).Value = 'Daniel 7:11'
$ws.Cells.Item(177, 1).Value = 'Daniel 7:12'
$ws.Cells.Item(178, 1).Value = 'Daniel 7:13'
$ws.Cells.Item(179, 1).Value = 'Daniel 7:14'
$ws.Cells.Item(180, 1).Value = 'Daniel 7:15'
$ws.Cells.Item(181, 1).Value = 'Daniel 7:16'
$ws.Cells.Item(182, 1).Value = 'Daniel 7:18'
$ws.Cells.Item(183, 1).Value = 'Daniel 7:19'
$ws.Cells.Item(184, 1).Value = 'Daniel 7:21'
$ws.Cells.Item(185, 1).Value = 'Daniel 7:22'
$ws.Cells.Item(186, 1).Value = 'Daniel 7:23'
$ws.Cells.Item(187, 1).Value = 'Daniel 7:24'
$ws.Cells.Item(188, 1).Value = 'Daniel 7:25'
$ws.Cells.Item(189, 1).Value = 'Daniel 7:26'
$ws.Cells.Item(190, 1).Value = 'Daniel 7:27'
$ws.Cells.Item(191, 1).Value = 'Daniel 7:28'
$ws.Cells.Item(192, 1).Value = 'Daniel 8:1'
$ws.Cells.Item(193, 1).Value = 'Daniel 8:2'
$ws.Cells.Item(194, 1).Value = 'Daniel 8:3'
$ws.Cells.Item(195, 1).Value = 'Daniel 8:4'
$ws.Cells.Item(196, 1).Value = 'Daniel 8:5'
$ws.Cells.Item(197, 1).Value = 'Daniel 8:6'
$ws.Cells.Item(198, 1).Value = 'Daniel 8:7'
$ws.Cells.Item(199, 1).Value = 'Daniel 8:8'
$ws.Cells.Item(200, 1).Value = 'Daniel 8:9'
$ws.Cells.Item(201, 1).Value = 'Daniel 8:12'
$ws.Cells.Item(202, 1).Value = 'Daniel 8:13'
$ws.Cells.Item(203, 1).Value = 'Daniel 8:14'
$ws.Cells.Item(204, 1).Value = 'Daniel 8:15'
$ws.Cells.Item(205, 1).Value = 'Daniel 8:17'
$ws.Cells.Item(206, 1).Value = 'Daniel 8:18'
$ws.Cells.Item(207, 1).Value = 'Daniel 8:19'
$ws.Cells.Item(208, 1).Value = 'Daniel 8:20'
$ws.Cells.Item(209, 1).Value = 'Daniel 8:21'
$ws.Cells.Item(210, 1).Value = 'Daniel 8:23'
$ws.Cells.Item(211, 1).Value = 'Daniel 8:24'
$ws.Cells.Item(212, 1).Value = 'Daniel 8:26'
$ws.Cells.Item(213, 1).Value = 'Daniel 8:27'
$ws.Cells.Item(214, 1).Value = 'Daniel 9:2'
$ws.Cells.Item(215, 1).Value = 'Daniel 9:3'
$ws.Cells.Item(216, 1).Value = 'Daniel 9:4'
$ws.Cells.Item(217, 1).Value = 'Daniel 9:7'
$ws.Cells.Item(218, 1).Value = 'Daniel 9:8'
$ws.Cells.Item(219, 1).Value = 'Daniel 9:12'
$ws.Cells.Item(220, 1).Value = 'Daniel 9:15'
$ws.Cells.Item(221, 1).Value = 'Daniel 9:16'
$ws.Cells.Item(222, 1).Value = 'Daniel 9:17'
$ws.Cells.Item(223, 1).Value = 'Daniel 9:19'
$ws.Cells.Item(224, 1).Value = 'Daniel 9:21'
$ws.Cells.Item(225, 1).Value = 'Daniel 9:26'
$ws.Cells.Item(226, 1).Value = 'Daniel 9:27'
$ws.Cells.Item(227, 1).Value = 'Daniel 10:2'
$ws.Cells.Item(228, 1).Value = 'Daniel 10:3'
$ws.Cells.Item(229, 1).Value = 'Daniel 10:4'
$ws.Cells.Item(230, 1).Value = 'Daniel 10:6'
$ws.Cells.Item(231, 1).Value = 'Daniel 10:8'
$ws.Cells.Item(232, 1).Value = 'Daniel 10:9'
$ws.Cells.Item(233, 1).Value = 'Daniel 10:10'
$ws.Cells.Item(234, 1).Value = 'Daniel 10:11'
$ws.Cells.Item(235, 1).Value = 'Daniel 10:12'
$ws.Cells.Item(236, 1).Value = 'Daniel 10:15'
$ws.Cells.Item(237, 1).Value = 'Daniel 10:16'
$ws.Cells.Item(238, 1).Value = 'Daniel 10:17'
$ws.Cells.Item(239, 1).Value = 'Daniel 10:18'
$ws.Cells.Item(240, 1).Value = 'Daniel 10:19'
$ws.Cells.Item(241, 1).Value = 'Daniel 10:20'
$ws.Cells.Item(242, 1).Value = 'Daniel 10:21'
$ws.Cells.Item(243, 1).Value = 'Daniel 11:1'
$ws.Cells.Item(244, 1).Value = 'Daniel 11:2'
$ws.Cells.Item(245, 1).Value = 'Daniel 11:3'
$ws.Cells.Item(246, 1).Value = 'Daniel 11:4'
$ws.Cells.Item(247, 1).Value = 'Daniel 11:5'
$ws.Cells.Item(248, 1).Value = 'Daniel 11:6'
$ws.Cells.Item(249, 1).Value = 'Daniel 11:7'
$ws.Cells.Item(250, 1).Value = 'Daniel 11:8'
$ws.Cells.Item(251, 1).Value = 'Daniel 11:9'
$ws.Cells.Item(252, 1).Value = 'Daniel 11:10'
$ws.Cells.Item(253, 1).Value = 'Daniel 11:11'
$ws.Cells.Item(254, 1).Value = 'Daniel 11:12'
$ws.Cells.Item(255, 1).Value = 'Daniel 11:13'
$ws.Cells.Item(256, 1).Value = 'Daniel 11:14'
$ws.Cells.Item(257, 1).Value = 'Daniel 11:15'
$ws.Cells.Item(258, 1).Value = 'Daniel 11:16'
$ws.Cells.Item(259, 1).Value = 'Daniel 11:17'
$ws.Cells.Item(260, 1).Value = 'Daniel 11:18'
$ws.Cells.Item(261, 1).Value = 'Daniel 11:19'
$ws.Cells.Item(262, 1).Value = 'Daniel 11:20'
$ws.Cells.Item(263, 1).Value = 'Daniel 11:21'
$ws.Cells.Item(264, 1).Value = 'Daniel 11:22'
$ws.Cells.Item(265, 1).Value = 'Daniel 11:23'
$ws.Cells.Item(266, 1).Value = 'Daniel 11:25'
$ws.Cells.Item(267, 1).Value = 'Daniel 11:26'
$ws.Cells.Item(268, 1).Value = 'Daniel 11:27'
$ws.Cells.Item(269, 1).Value = 'Daniel 11:28'
$ws.Cells.Item(270, 1).Value = 'Daniel 11:29'
$ws.Cells.Item(271, 1).Value = 'Daniel 11:30'
$ws.Cells.Item(272, 1).Value = 'Daniel 11:31'
$ws.Cells.Item(273, 1).Value = 'Daniel 11:32'
$ws.Cells.Item(274, 1).Value = 'Daniel 11:34'
$ws.Cells.Item(275, 1).Value = 'Daniel 11:35'
$ws.Cells.Item(276, 1).Value = 'Daniel 11:36'
$ws.Cells.Item(277, 1).Value = 'Daniel 11:37'
$ws.Cells.Item(278, 1).Value = 'Daniel 11:38'
$ws.Cells.Item(279, 1).Value = 'Daniel 11:39'
$ws.Cells.Item(280, 1).Value = 'Daniel 11:41'
$ws.Cells.Item(281, 1).Value = 'Daniel 11:42'
$ws.Cells.Item(282, 1).Value = 'Daniel 11:43'
$ws.Cells.Item(283, 1).Value = 'Daniel 11:44'
$ws.Cells.Item(284, 1).Value = 'Daniel 11:45'
$ws.Cells.Item(285, 1).Value = 'Daniel 12:1'
$ws.Cells.Item(286, 1).Value = 'Daniel 12:2'
$ws.Cells.Item(287, 1).Value = 'Daniel 12:4'
$ws.Cells.Item(288, 1).Value = 'Daniel 12:5'
$ws.Cells.Item(289, 1).Value = 'Daniel 12:6'
$ws.Cells.Item(290, 1).Value = 'Daniel 12:7'
$ws.Cells.Item(291, 1).Value = 'Daniel 12:8'
$ws.Cells.Item(292, 1).Value = 'Daniel 12:9'
$ws.Cells.Item(293, 1).Value = 'Daniel 12:11'
